$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 10056
$ws.Range("B3").Value = 10000
$ws.Range("C3").Value = 316.81
$ws.Range("D3").Value = 318.58999999999997
$ws.Range("E3").Value = $false
$ws.Range("F3").Value = 0.56000000000000005
$ws.Range("G3").Value = 42606.42832175926
$ws.Range("H3").Value = $true

# Row 4
$ws.Range("A4").Value = 10125.39
$ws.Range("B4").Value = 10056
$ws.Range("C4").Value = 316.81
$ws.Range("D4").Value = 319
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.69
$ws.Range("G4").Value = 42606.48678240741
$ws.Range("H4").Value = $true

$ws.Columns.Item(1).ColumnWidth = 8.14
